$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row before row 4. This shifts the old row 4 ("Number of
#    disability persons" data) down to row 5, and the old row 5 (merged
#    source citation) down to row 6, making room for the new
#    "family with disabilities Persons" data row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. Row 1 - new title, merged across A1:I1
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tkibuli Municipality"
$ws.Rows.Item(1).RowHeight = 51
$ws.Range("A1:I1").Font.Name = "Arial"
$ws.Range("A1:I1").Font.Size = 11
$ws.Range("A1:I1").Font.Bold = $true
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row 2 - unchanged text, but no longer has a custom row height
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------------------
# 4. Row 3 - A3 font changes to Sylfaen; years stay the same
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 18.75
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------------------
# 5. Row 4 (new) - "family with disabilities Persons"
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 24.75

$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(9).LineStyle = 0

$row4vals = @(619,595,534,524,513,511,496,499)
$cols = @("B","C","D","E","F","G","H","I")
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "4")
    $cell.Value = $row4vals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.NumberFormat = "#\ ##0"
    $cell.Borders.Item(8).LineStyle = 0
    $cell.Borders.Item(9).LineStyle = 0
}

# ---------------------------------------------------------------------------
# 6. Row 5 (was old row 4) - "disabilities Persons"
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 21

$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1

$row5vals = @(677,649,583,567,556,553,538,544)
for ($i = 0; $i -lt 8; $i++) {
    $cell = $ws.Range($cols[$i] + "5")
    $cell.Value = $row5vals[$i]
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.NumberFormat = "#\ ##0"
    $cell.Borders.Item(8).LineStyle = 0
    $cell.Borders.Item(9).LineStyle = 0
}
# I5 keeps a bottom border (matches the bottom of the data block)
$ws.Range("I5").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 7. Row 6 (was old row 5, the merged source citation) - row height only
#    (text/merge/rich-run formatting already carried forward by the insert)
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------------------
# 8. Column widths - column A widens, others revert to the workbook default
# ---------------------------------------------------------------------------
$ws.Range("B1:P1").EntireColumn.ColumnWidth = 8.25
$ws.Columns.Item(1).ColumnWidth = 19.92
